$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two leading blank rows so the header (old row 3) becomes row 1
# and the data block shifts up to match the new A1:G65 layout.
$ws.Rows("1:2").Delete()

# Re-apply the AutoFilter over the new extent of the table.
$ws.AutoFilterMode = $false
$ws.Range("A1:G65").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the
# new AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=TABLA!`$A`$1:`$G`$65"
    }
}

# Restore the selected cell to match the saved view.
$ws.Range("J8").Select()
